$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H15").Value = 219020.97
$ws_ALC.Range("I15").Value = 219020.97
$ws_ALC.Range("K15").Value = 657062.91
$ws_ALC.Range("M15").Value = -656893.91
$ws_ALC.Range("H106").Value = 1663.3334
$ws_ALC.Range("I106").Value = 1663.3334
$ws_ALC.Range("J106").Value = 0
$ws_ALC.Range("K106").Value = 1663.3334
$ws_ALC.Range("L106").Value = 0
$ws_ALC.Range("N106").ClearContents()
$ws_ALC.Range("H107").Value = 481476.9
$ws_ALC.Range("I107").Value = 721681.8
$ws_ALC.Range("J107").Value = 1067.1428
$ws_ALC.Range("K107").Value = 721681.8
$ws_ALC.Range("L107").Value = 1067.1428
$ws_ALC.Range("M107").Value = -719761.8
$ws_ALC.Range("N107").Value = -4907.1428
$ws_ALC.Range("H113").Value = 13573.75
$ws_ALC.Range("I113").Value = 3905
$ws_ALC.Range("K113").Value = 3905
$ws_ALC.Range("M113").Value = -651
$ws_ALC.Range("H116").Value = 2685
$ws_ALC.Range("I116").Value = 2714.2856
$ws_ALC.Range("J116").Value = 2480
$ws_ALC.Range("K116").Value = 2714.2856
$ws_ALC.Range("L116").Value = 2480
$ws_ALC.Range("M116").Value = 727.7143999999998
$ws_ALC.Range("N116").Value = -9364
$ws_ALC.Range("H132").Value = 2999.6
$ws_ALC.Range("I132").Value = 3067.7273
$ws_ALC.Range("J132").Value = 2500
$ws_ALC.Range("K132").Value = 9203.1819
$ws_ALC.Range("L132").Value = 7500
$ws_ALC.Range("M132").Value = -6673.1819
$ws_ALC.Range("N132").Value = -12560
$ws_ALC.Range("H137").Value = 1143.6471
$ws_ALC.Range("I137").Value = 1083.8334
$ws_ALC.Range("J137").Value = 1422.7778
$ws_ALC.Range("K137").Value = 3251.5002
$ws_ALC.Range("L137").Value = 4268.3334
$ws_ALC.Range("M137").Value = -701.5001999999999
$ws_ALC.Range("N137").Value = -9368.3334
$ws_ALC.Range("H141").Value = 5092.033
$ws_ALC.Range("I141").Value = 1082.6459
$ws_ALC.Range("J141").Value = 21129.584
$ws_ALC.Range("K141").Value = 3247.9377
$ws_ALC.Range("L141").Value = 63388.75199999999
$ws_ALC.Range("M141").Value = 1932.0623
$ws_ALC.Range("N141").Value = -73748.75199999999
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 4132.1157
$ws_ARM.Range("I32").Value = 3106.055
$ws_ARM.Range("J32").Value = 27475
$ws_ARM.Range("K32").Value = 3106.055
$ws_ARM.Range("L32").Value = 27475
$ws_ARM.Range("M32").Value = -2819.055
$ws_ARM.Range("N32").Value = -28049
$ws_ARM.Range("H61").Value = 766.175
$ws_ARM.Range("I61").Value = 573.21875
$ws_ARM.Range("J61").Value = 1538
$ws_ARM.Range("K61").Value = 573.21875
$ws_ARM.Range("L61").Value = 1538
$ws_ARM.Range("M61").Value = -361.21875
$ws_ARM.Range("N61").Value = -1962
$ws_ARM.Range("H74").Value = 3559.682
$ws_ARM.Range("I74").Value = 3851.7437
$ws_ARM.Range("J74").Value = 1281.6
$ws_ARM.Range("K74").Value = 3851.7437
$ws_ARM.Range("L74").Value = 1281.6
$ws_ARM.Range("M74").Value = -2977.7437
$ws_ARM.Range("N74").Value = -3029.6
$ws_ARM.Range("H77").Value = 3559.682
$ws_ARM.Range("I77").Value = 3851.7437
$ws_ARM.Range("J77").Value = 1281.6
$ws_ARM.Range("K77").Value = 19258.7185
$ws_ARM.Range("L77").Value = 6408
$ws_ARM.Range("M77").Value = -14890.7185
$ws_ARM.Range("N77").Value = -15144
$ws_ARM.Range("H132").Value = 2739.4849
$ws_ARM.Range("I132").Value = 1486.7858
$ws_ARM.Range("J132").Value = 3662.5264
$ws_ARM.Range("K132").Value = 4460.357400000001
$ws_ARM.Range("L132").Value = 10987.5792
$ws_ARM.Range("M132").Value = -1930.357400000001
$ws_ARM.Range("N132").Value = -16047.5792
$ws_ARM.Range("H136").Value = 766.175
$ws_ARM.Range("I136").Value = 573.21875
$ws_ARM.Range("J136").Value = 1538
$ws_ARM.Range("K136").Value = 1719.65625
$ws_ARM.Range("L136").Value = 4614
$ws_ARM.Range("M136").Value = 830.34375
$ws_ARM.Range("N136").Value = -9714
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H97").Value = 17366.334
$ws_BSM.Range("J97").Value = 25735.5
$ws_BSM.Range("L97").Value = 25735.5
$ws_BSM.Range("N97").Value = -27717.5
$ws_BSM.Range("H99").Value = 1485
$ws_BSM.Range("I99").Value = 1703.3334
$ws_BSM.Range("K99").Value = 1703.3334
$ws_BSM.Range("M99").Value = -205.3334
$ws_BSM.Range("H134").Value = 2060.1667
$ws_BSM.Range("I134").Value = 1215.25
$ws_BSM.Range("J134").Value = 3750
$ws_BSM.Range("K134").Value = 3645.75
$ws_BSM.Range("L134").Value = 11250
$ws_BSM.Range("M134").Value = -1110.75
$ws_BSM.Range("N134").Value = -16320
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 2594.9355
$ws_CRP.Range("I31").Value = 1754.5769
$ws_CRP.Range("J31").Value = 3201.861
$ws_CRP.Range("K31").Value = 1754.5769
$ws_CRP.Range("L31").Value = 3201.861
$ws_CRP.Range("M31").Value = -1459.5769
$ws_CRP.Range("N31").Value = -3791.861
$ws_CRP.Range("H34").Value = 2594.9355
$ws_CRP.Range("I34").Value = 1754.5769
$ws_CRP.Range("J34").Value = 3201.861
$ws_CRP.Range("K34").Value = 1754.5769
$ws_CRP.Range("L34").Value = 3201.861
$ws_CRP.Range("M34").Value = -1552.5769
$ws_CRP.Range("N34").Value = -3605.861
$ws_CRP.Range("H58").Value = 1079.403
$ws_CRP.Range("I58").Value = 924.2712
$ws_CRP.Range("J58").Value = 2223.5
$ws_CRP.Range("K58").Value = 924.2712
$ws_CRP.Range("L58").Value = 2223.5
$ws_CRP.Range("M58").Value = -721.2712
$ws_CRP.Range("N58").Value = -2629.5
$ws_CRP.Range("H99").Value = 2920100
$ws_CRP.Range("I99").Value = 6401760
$ws_CRP.Range("J99").Value = 18716.666
$ws_CRP.Range("K99").Value = 6401760
$ws_CRP.Range("L99").Value = 18716.666
$ws_CRP.Range("M99").Value = -6400262
$ws_CRP.Range("N99").Value = -21712.666
$ws_CRP.Range("H126").Value = 2920100
$ws_CRP.Range("I126").Value = 6401760
$ws_CRP.Range("J126").Value = 18716.666
$ws_CRP.Range("K126").Value = 19205280
$ws_CRP.Range("L126").Value = 56149.99800000001
$ws_CRP.Range("M126").Value = -19202810
$ws_CRP.Range("N126").Value = -61089.99800000001
$ws_CRP.Range("H132").Value = 2468.2173
$ws_CRP.Range("I132").Value = 1212.1666
$ws_CRP.Range("J132").Value = 3838.4546
$ws_CRP.Range("K132").Value = 3636.4998
$ws_CRP.Range("L132").Value = 11515.3638
$ws_CRP.Range("M132").Value = -1106.4998
$ws_CRP.Range("N132").Value = -16575.3638
$ws_CRP.Range("H134").Value = 1835.9688
$ws_CRP.Range("I134").Value = 1997.7307
$ws_CRP.Range("J134").Value = 1135
$ws_CRP.Range("K134").Value = 5993.1921
$ws_CRP.Range("L134").Value = 3405
$ws_CRP.Range("M134").Value = -3458.1921
$ws_CRP.Range("N134").Value = -8475
$ws_CRP.Range("H136").Value = 1079.403
$ws_CRP.Range("I136").Value = 924.2712
$ws_CRP.Range("J136").Value = 2223.5
$ws_CRP.Range("K136").Value = 2772.8136
$ws_CRP.Range("L136").Value = 6670.5
$ws_CRP.Range("M136").Value = -222.8136
$ws_CRP.Range("N136").Value = -11770.5
$ws_CRP.Range("H138").Value = 64972.5
$ws_CRP.Range("J138").Value = 64972.5
$ws_CRP.Range("L138").Value = 64972.5
$ws_CRP.Range("N138").Value = -75252.5
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H122").Value = 782.4091
$ws_CUL.Range("J122").Value = 1123.909
$ws_CUL.Range("L122").Value = 10115.181
$ws_CUL.Range("N122").Value = -15015.181
$ws_CUL.Range("H132").Value = 1438
$ws_CUL.Range("I132").Value = 1577.8889
$ws_CUL.Range("J132").Value = 1378.0476
$ws_CUL.Range("K132").Value = 14201.0001
$ws_CUL.Range("L132").Value = 12402.4284
$ws_CUL.Range("M132").Value = -11671.0001
$ws_CUL.Range("N132").Value = -17462.4284
$ws_CUL.Range("H137").Value = 2077.2144
$ws_CUL.Range("I137").Value = 2266.125
$ws_CUL.Range("J137").Value = 2001.65
$ws_CUL.Range("K137").Value = 6798.375
$ws_CUL.Range("L137").Value = 6004.950000000001
$ws_CUL.Range("M137").Value = -1698.375
$ws_CUL.Range("N137").Value = -16204.95
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H132").Value = 2349.8
$ws_GSM.Range("I132").Value = 2082.1765
$ws_GSM.Range("J132").Value = 2699.7693
$ws_GSM.Range("K132").Value = 6246.529500000001
$ws_GSM.Range("L132").Value = 8099.3079
$ws_GSM.Range("M132").Value = -3716.529500000001
$ws_GSM.Range("N132").Value = -13159.3079
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H122").Value = 1652.8
$ws_LTW.Range("I122").Value = 1266.1428
$ws_LTW.Range("J122").Value = 2555
$ws_LTW.Range("K122").Value = 3798.4284
$ws_LTW.Range("L122").Value = 7665
$ws_LTW.Range("M122").Value = -1348.4284
$ws_LTW.Range("N122").Value = -12565
$ws_LTW.Range("H132").Value = 5794.8335
$ws_LTW.Range("I132").Value = 5819.6807
$ws_LTW.Range("J132").Value = 5748.12
$ws_LTW.Range("K132").Value = 17459.0421
$ws_LTW.Range("L132").Value = 17244.36
$ws_LTW.Range("M132").Value = -14929.0421
$ws_LTW.Range("N132").Value = -22304.36
$ws_LTW.Range("H136").Value = 2827.7856
$ws_LTW.Range("I136").Value = 2827.7856
$ws_LTW.Range("K136").Value = 8483.356800000001
$ws_LTW.Range("M136").Value = -5933.356800000001
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H132").Value = 1417.6666
$ws_WVR.Range("I132").Value = 1388.4762
$ws_WVR.Range("J132").Value = 1519.8334
$ws_WVR.Range("K132").Value = 4165.4286
$ws_WVR.Range("L132").Value = 4559.5002
$ws_WVR.Range("M132").Value = -1635.4286
$ws_WVR.Range("N132").Value = -9619.5002
$ws_WVR.Range("H136").Value = 1688.2325
$ws_WVR.Range("I136").Value = 624.56757
$ws_WVR.Range("J136").Value = 8247.5
$ws_WVR.Range("K136").Value = 1873.70271
$ws_WVR.Range("L136").Value = 24742.5
$ws_WVR.Range("M136").Value = 676.29729
$ws_WVR.Range("N136").Value = -29842.5
$ws_WVR.Range("H139").Value = 65435
$ws_WVR.Range("I139").Value = 40925
$ws_WVR.Range("J139").Value = 71562.5
$ws_WVR.Range("K139").Value = 40925
$ws_WVR.Range("L139").Value = 71562.5
$ws_WVR.Range("M139").Value = -35785
$ws_WVR.Range("N139").Value = -81842.5
